$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 15-column x 2-row dump and rebuild as a 2-column Key/Value table
$ws.Cells.Clear()

$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Value"

$ws.Range("A2").Value = "RemoteAccess"
$ws.Range("B2").Value = "yes"

$ws.Range("A3").Value = "Site table"
$ws.Range("B3").Value = "wee"

$ws.Range("A4").Value = "Device"
$ws.Range("B4").Value = "abc"

$ws.Range("A5").Value = "FriendlyName"
$ws.Range("B5").Value = "sxk"

$ws.Range("A6").Value = "OS"
$ws.Range("B6").Value = "sdff"

$ws.Range("A7").Value = "NOCDefined"
$ws.Range("B7").Value = "dfasdf"

# Widen the key column like the target sheet
$ws.Columns.Item(1).ColumnWidth = 23.67

# Leave the selection on the last cell, matching the saved view state
$ws.Range("B7").Select()
